$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (Fecha, Calidad, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Precio $/Kg) got reshuffled across rows as
# part of the weekly update. Apply the new values directly per row.

# Row 5
$ws.Range("D5").Value = 44224
$ws.Range("K5").Value = 850
$ws.Range("L5").Value = 900
$ws.Range("M5").Value = 875
$ws.Range("P5").Value = 875

# Row 6
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 800
$ws.Range("K6").Value = 450
$ws.Range("L6").Value = 500
$ws.Range("M6").Value = 475
$ws.Range("P6").Value = 475

# Row 7
$ws.Range("D7").Value = 44174
$ws.Range("I7").Value = "Tercera"
$ws.Range("J7").Value = 1200
$ws.Range("K7").Value = 250
$ws.Range("L7").Value = 350
$ws.Range("M7").Value = 300
$ws.Range("P7").Value = 300

# Row 8
$ws.Range("D8").Value = 44278
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 700
$ws.Range("K8").Value = 600
$ws.Range("L8").Value = 700
$ws.Range("M8").Value = 650
$ws.Range("P8").Value = 650

# Row 9
$ws.Range("D9").Value = 44278
$ws.Range("I9").Value = "Tercera"
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 500
$ws.Range("L9").Value = 600
$ws.Range("M9").Value = 550
$ws.Range("P9").Value = 550

# Row 10
$ws.Range("D10").Value = 44245
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 800
$ws.Range("K10").Value = 850
$ws.Range("L10").Value = 900
$ws.Range("M10").Value = 875
$ws.Range("P10").Value = 875

# Row 11
$ws.Range("D11").Value = 44245
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 750
$ws.Range("L11").Value = 800
$ws.Range("M11").Value = 775
$ws.Range("P11").Value = 775

# Row 12
$ws.Range("D12").Value = 44229
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 760
$ws.Range("K12").Value = 550
$ws.Range("M12").Value = 575
$ws.Range("P12").Value = 575

# Row 14
$ws.Range("D14").Value = 44210
$ws.Range("J14").Value = 900
$ws.Range("K14").Value = 600
$ws.Range("L14").Value = 700
$ws.Range("M14").Value = 650
$ws.Range("P14").Value = 650
